$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 00:19"

# Update country statistics rows (col order: A=Pais, B=Casos totales, C=Nuevos casos,
# D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 7349869
$ws.Cells.Item(4,3).Value = 28526
$ws.Cells.Item(4,4).Value = 4599310
$ws.Cells.Item(4,5).Value = 2540850
$ws.Cells.Item(4,7).Value = 256
$ws.Cells.Item(4,8).Value = 209709

# Row 29: Canada
$ws.Cells.Item(29,2).Value = 155034
$ws.Cells.Item(29,3).Value = 1909
$ws.Cells.Item(29,4).Value = 132297
$ws.Cells.Item(29,5).Value = 13462
$ws.Cells.Item(29,7).Value = 7
$ws.Cells.Item(29,8).Value = 9275

# Row 41: Egipto
$ws.Cells.Item(41,2).Value = 102955
$ws.Cells.Item(41,3).Value = 115
$ws.Cells.Item(41,4).Value = 95586
$ws.Cells.Item(41,5).Value = 1468
$ws.Cells.Item(41,7).Value = 18
$ws.Cells.Item(41,8).Value = 5901

# Row 48: Japon
$ws.Cells.Item(48,2).Value = 82131
$ws.Cells.Item(48,3).Value = 441
$ws.Cells.Item(48,4).Value = 74990
$ws.Cells.Item(48,5).Value = 5593
$ws.Cells.Item(48,7).Value = 3
$ws.Cells.Item(48,8).Value = 1548

# Row 56: Barein
$ws.Cells.Item(56,2).Value = 69848
$ws.Cells.Item(56,3).Value = 487
$ws.Cells.Item(56,4).Value = 63549
$ws.Cells.Item(56,5).Value = 6054

# Row 65: Ghana
$ws.Cells.Item(65,2).Value = 46444
$ws.Cells.Item(65,3).Value = 57
$ws.Cells.Item(65,4).Value = 45646
$ws.Cells.Item(65,5).Value = 499

# Row 83: Camerun
$ws.Cells.Item(83,2).Value = 20838
$ws.Cells.Item(83,3).Value = 103
$ws.Cells.Item(83,5).Value = 980

# Row 84: Bulgaria
$ws.Cells.Item(84,2).Value = 20271
$ws.Cells.Item(84,3).Value = 216
$ws.Cells.Item(84,4).Value = 14339
$ws.Cells.Item(84,5).Value = 5125
$ws.Cells.Item(84,7).Value = 11
$ws.Cells.Item(84,8).Value = 807

# Row 108: Haiti
$ws.Cells.Item(108,4).Value = 6757
$ws.Cells.Item(108,5).Value = 1756

# Row 112: Zimbabue
$ws.Cells.Item(112,2).Value = 7816
$ws.Cells.Item(112,3).Value = 4
$ws.Cells.Item(112,4).Value = 6112
$ws.Cells.Item(112,5).Value = 1476
$ws.Cells.Item(112,7).Value = 1
$ws.Cells.Item(112,8).Value = 228

# Row 136: Aruba
$ws.Cells.Item(136,2).Value = 3872
$ws.Cells.Item(136,3).Value = 28
$ws.Cells.Item(136,4).Value = 3037
$ws.Cells.Item(136,5).Value = 810

# Rows 144-146: Botsuana overtakes Mali and Malta in ranking (re-sort by Casos totales)
$ws.Cells.Item(144,1).Value = "Botsuana"
$ws.Cells.Item(144,2).Value = 3172
$ws.Cells.Item(144,3).Value = 251
$ws.Cells.Item(144,4).Value = 710
$ws.Cells.Item(144,5).Value = 2446
$ws.Cells.Item(144,7).Value = 0
$ws.Cells.Item(144,8).Value = 16

$ws.Cells.Item(145,1).Value = "Mali"
$ws.Cells.Item(145,2).Value = 3090
$ws.Cells.Item(145,3).Value = 4
$ws.Cells.Item(145,4).Value = 2439
$ws.Cells.Item(145,5).Value = 521
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 130

$ws.Cells.Item(146,1).Value = "Malta"
$ws.Cells.Item(146,2).Value = 3006
$ws.Cells.Item(146,3).Value = 27
$ws.Cells.Item(146,4).Value = 2399
$ws.Cells.Item(146,5).Value = 575
$ws.Cells.Item(146,7).Value = 1
$ws.Cells.Item(146,8).Value = 32

# Row 166: Niger
$ws.Cells.Item(166,4).Value = 1110
$ws.Cells.Item(166,5).Value = 17

# Rows 207-208: Santa Lucia and Timor Oriental swap order (tied case counts)
$ws.Cells.Item(207,1).Value = "Santa Lucia"
$ws.Cells.Item(208,1).Value = "Timor Oriental"
